$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new task ("Fix normals problem on cylinder") was added to the top of the
# list, pushing the rest of the rows down by one. Insert a fresh row at 2
# (this is the row that currently holds "Data drive lighting").
$ws.Rows.Item(2).Insert()

# Row insert copies the header row's (bold) formatting down onto the new
# row - put it back to normal/unbold so it matches the rest of the table.
$ws.Range("A2:B2").Font.Bold = $false

# Fill in the new row with the new task and its estimate.
$ws.Range("A2").Value = "Fix normals problem on cylinder"
$ws.Range("B2").Value = 3

# "Data drive lighting" (now shifted down to row 3) gets a revised estimate.
$ws.Range("B3").Value = 4

# Match the workbook's on-disk selection state.
$ws.Range("A2:B2").Select()
